$wb = $excel.ActiveWorkbook

# Data refresh: update market price / profit figures across sheets
# as produced by the scheduled runner (raw values, no formulas).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3316.5518
$ws.Range("J113").Value = 2854.4443
$ws.Range("L113").Value = 2854.4443
$ws.Range("N113").Value = -9362.444299999999

$ws.Range("H116").Value = 5034.625
$ws.Range("I116").Value = 5231.154
$ws.Range("K116").Value = 5231.154
$ws.Range("M116").Value = -1789.154

$ws.Range("H132").Value = 4002857.8
$ws.Range("I132").Value = 5002833
$ws.Range("K132").Value = 15008499
$ws.Range("M132").Value = -15005969

$ws.Range("H137").Value = 3114.0193
$ws.Range("I137").Value = 3060.9556
$ws.Range("J137").Value = 3455.1428
$ws.Range("K137").Value = 9182.8668
$ws.Range("L137").Value = 10365.4284
$ws.Range("M137").Value = -6632.8668
$ws.Range("N137").Value = -15465.4284

$ws.Range("H138").Value = 2515.9878
$ws.Range("I138").Value = 1010.5106
$ws.Range("J138").Value = 4537.6284
$ws.Range("K138").Value = 3031.5318
$ws.Range("L138").Value = 13612.8852
$ws.Range("M138").Value = 2108.4682
$ws.Range("N138").Value = -23892.8852

$ws.Range("H141").Value = 155764.23
$ws.Range("I141").Value = 2611.242
$ws.Range("J141").Value = 1738345.1
$ws.Range("K141").Value = 7833.726000000001
$ws.Range("L141").Value = 5215035.300000001
$ws.Range("M141").Value = -2653.726000000001
$ws.Range("N141").Value = -5225395.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1758.19
$ws.Range("I32").Value = 1503.5483
$ws.Range("J32").Value = 5141.2856
$ws.Range("K32").Value = 1503.5483
$ws.Range("L32").Value = 5141.2856
$ws.Range("M32").Value = -1216.5483
$ws.Range("N32").Value = -5715.2856

$ws.Range("H45").Value = 1395.0264
$ws.Range("I45").Value = 1081.5186
$ws.Range("J45").Value = 2164.5454
$ws.Range("K45").Value = 1081.5186
$ws.Range("L45").Value = 2164.5454
$ws.Range("M45").Value = -704.5186000000001
$ws.Range("N45").Value = -2918.5454

$ws.Range("H61").Value = 1230.3062
$ws.Range("I61").Value = 679.439
$ws.Range("J61").Value = 4053.5
$ws.Range("K61").Value = 679.439
$ws.Range("L61").Value = 4053.5
$ws.Range("M61").Value = -467.439
$ws.Range("N61").Value = -4477.5

$ws.Range("H74").Value = 628.2889
$ws.Range("I74").Value = 483.58975
$ws.Range("J74").Value = 1568.8334
$ws.Range("K74").Value = 483.58975
$ws.Range("L74").Value = 1568.8334
$ws.Range("M74").Value = 390.41025
$ws.Range("N74").Value = -3316.8334

$ws.Range("H77").Value = 628.2889
$ws.Range("I77").Value = 483.58975
$ws.Range("J77").Value = 1568.8334
$ws.Range("K77").Value = 2417.94875
$ws.Range("L77").Value = 7844.166999999999
$ws.Range("M77").Value = 1950.05125
$ws.Range("N77").Value = -16580.167

$ws.Range("H102").Value = 3032.0667
$ws.Range("I102").Value = 2498.5386
$ws.Range("K102").Value = 2498.5386
$ws.Range("M102").Value = -876.5385999999999

$ws.Range("H122").Value = 3435.64
$ws.Range("I122").Value = 2894.55
$ws.Range("K122").Value = 8683.650000000001
$ws.Range("M122").Value = -6233.650000000001

$ws.Range("H132").Value = 1648.0492
$ws.Range("I132").Value = 1246.674
$ws.Range("J132").Value = 2878.9333
$ws.Range("K132").Value = 3740.022
$ws.Range("L132").Value = 8636.7999
$ws.Range("M132").Value = -1210.022
$ws.Range("N132").Value = -13696.7999

$ws.Range("H136").Value = 1230.3062
$ws.Range("I136").Value = 679.439
$ws.Range("J136").Value = 4053.5
$ws.Range("K136").Value = 2038.317
$ws.Range("L136").Value = 12160.5
$ws.Range("M136").Value = 511.683
$ws.Range("N136").Value = -17260.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

$ws.Range("H134").Value = 2833.6924
$ws.Range("I134").Value = 2707.652
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 8122.956
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -5587.956
$ws.Range("N134").Value = -16470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2677.151
$ws.Range("I31").Value = 1502.2821
$ws.Range("J31").Value = 5950
$ws.Range("K31").Value = 1502.2821
$ws.Range("L31").Value = 5950
$ws.Range("M31").Value = -1207.2821
$ws.Range("N31").Value = -6540

$ws.Range("H34").Value = 2677.151
$ws.Range("I34").Value = 1502.2821
$ws.Range("J34").Value = 5950
$ws.Range("K34").Value = 1502.2821
$ws.Range("L34").Value = 5950
$ws.Range("M34").Value = -1300.2821
$ws.Range("N34").Value = -6354

$ws.Range("H58").Value = 7814768.5
$ws.Range("I58").Value = 1493.06
$ws.Range("J58").Value = 35719324
$ws.Range("K58").Value = 1493.06
$ws.Range("L58").Value = 35719324
$ws.Range("M58").Value = -1290.06
$ws.Range("N58").Value = -35719730

$ws.Range("H132").Value = 1838.25
$ws.Range("I132").Value = 1419.2162
$ws.Range("J132").Value = 2871.8667
$ws.Range("K132").Value = 4257.6486
$ws.Range("L132").Value = 8615.6001
$ws.Range("M132").Value = -1727.6486
$ws.Range("N132").Value = -13675.6001

$ws.Range("H134").Value = 1151.847
$ws.Range("I134").Value = 820.7361
$ws.Range("J134").Value = 2985.6924
$ws.Range("K134").Value = 2462.2083
$ws.Range("L134").Value = 8957.0772
$ws.Range("M134").Value = 72.79170000000022
$ws.Range("N134").Value = -14027.0772

$ws.Range("H136").Value = 7814768.5
$ws.Range("I136").Value = 1493.06
$ws.Range("J136").Value = 35719324
$ws.Range("K136").Value = 4479.18
$ws.Range("L136").Value = 107157972
$ws.Range("M136").Value = -1929.18
$ws.Range("N136").Value = -107163072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7197.857
$ws.Range("I87").Value = 2958.077
$ws.Range("J87").Value = 14087.5
$ws.Range("K87").Value = 8874.231
$ws.Range("L87").Value = 42262.5
$ws.Range("M87").Value = -7626.231
$ws.Range("N87").Value = -44758.5

$ws.Range("H90").Value = 7197.857
$ws.Range("I90").Value = 2958.077
$ws.Range("J90").Value = 14087.5
$ws.Range("K90").Value = 26622.693
$ws.Range("L90").Value = 126787.5
$ws.Range("M90").Value = -20382.693
$ws.Range("N90").Value = -139267.5

$ws.Range("H113").Value = 1777.875
$ws.Range("I113").Value = 4321
$ws.Range("J113").Value = 930.1667
$ws.Range("K113").Value = 12963
$ws.Range("L113").Value = 2790.5001
$ws.Range("M113").Value = -10793
$ws.Range("N113").Value = -7130.5001

$ws.Range("H131").Value = 2164.96
$ws.Range("J131").Value = 1638.5555
$ws.Range("L131").Value = 4915.666499999999
$ws.Range("N131").Value = -14995.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 21292.46
$ws.Range("I102").Value = 1498.3823
$ws.Range("J102").Value = 58681.277
$ws.Range("K102").Value = 1498.3823
$ws.Range("L102").Value = 58681.277
$ws.Range("M102").Value = 123.6177
$ws.Range("N102").Value = -61925.277

$ws.Range("H122").Value = 3343.9688
$ws.Range("I122").Value = 2358.2632
$ws.Range("J122").Value = 4784.615
$ws.Range("K122").Value = 7074.7896
$ws.Range("L122").Value = 14353.845
$ws.Range("M122").Value = -4624.7896
$ws.Range("N122").Value = -19253.845

$ws.Range("H132").Value = 2144.9583
$ws.Range("I132").Value = 1871.5646
$ws.Range("K132").Value = 5614.6938
$ws.Range("M132").Value = -3084.6938

$ws.Range("H133").Value = 24000
$ws.Range("J133").Value = 24000
$ws.Range("L133").Value = 24000
$ws.Range("N133").Value = -34120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1622.3529
$ws.Range("I46").Value = 469.7143
$ws.Range("J46").Value = 2429.2
$ws.Range("K46").Value = 469.7143
$ws.Range("L46").Value = 2429.2
$ws.Range("M46").Value = -281.7143
$ws.Range("N46").Value = -2805.2

$ws.Range("H132").Value = 2343.8484
$ws.Range("I132").Value = 1321.3182
$ws.Range("J132").Value = 4388.909
$ws.Range("K132").Value = 3963.9546
$ws.Range("L132").Value = 13166.727
$ws.Range("M132").Value = -1433.9546
$ws.Range("N132").Value = -18226.727

$ws.Range("H136").Value = 1384.8806
$ws.Range("I136").Value = 944.62067
$ws.Range("J136").Value = 4222.1113
$ws.Range("K136").Value = 2833.86201
$ws.Range("L136").Value = 12666.3339
$ws.Range("M136").Value = -283.8620099999998
$ws.Range("N136").Value = -17766.3339

$ws.Range("H140").Value = 29076.334
$ws.Range("J140").Value = 29076.334
$ws.Range("L140").Value = 29076.334
$ws.Range("N140").Value = -39436.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 347238.7
$ws.Range("I122").Value = 501990.9
$ws.Range("J122").Value = 3344.889
$ws.Range("K122").Value = 1505972.7
$ws.Range("L122").Value = 10034.667
$ws.Range("M122").Value = -1503522.7
$ws.Range("N122").Value = -14934.667

$ws.Range("H132").Value = 2939.2622
$ws.Range("I132").Value = 898.525
$ws.Range("J132").Value = 6826.381
$ws.Range("K132").Value = 2695.575
$ws.Range("L132").Value = 20479.143
$ws.Range("M132").Value = -165.5749999999998
$ws.Range("N132").Value = -25539.143

$ws.Range("H136").Value = 985.3182
$ws.Range("I136").Value = 473.63333
$ws.Range("J136").Value = 2081.7856
$ws.Range("K136").Value = 1420.89999
$ws.Range("L136").Value = 6245.3568
$ws.Range("M136").Value = 1129.10001
$ws.Range("N136").Value = -11345.3568
